$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A (data currently ends at row 89)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

# Append two new training examples, each tagged as a "Вопрос" (question)
$ws.Cells.Item($newRow1, 1).Value = "мис   как   дело ?"
$ws.Cells.Item($newRow1, 2).Value = "Вопрос"
$ws.Cells.Item($newRow1, 3).Value = 1

$ws.Cells.Item($newRow2, 1).Value = "мис   как   погода ?"
$ws.Cells.Item($newRow2, 2).Value = "Вопрос"
$ws.Cells.Item($newRow2, 3).Value = 1
